$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "BlockMove - Animation could not found at assassine"
$ws.Range("B5").Value = "Not Fixed"

$ws.Range("A6").Value = "get stunned while in air --> slow fall"
$ws.Range("B6").Value = "Not Fixed"

$ws.Range("A9").Select()
